$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1986183074265976
$ws.Range("C2").Value = 0.5457685664939551
$ws.Range("J2").Value = 0.008635578583765112
$ws.Range("P2").Value = 0.150259067357513
$ws.Range("S2").Value = 0.09671848013816926
$ws.Range("B3").Value = 0.009375
$ws.Range("C3").Value = 0.021875
$ws.Range("J3").Value = 0.028125
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.190625
$ws.Range("J4").Value = 0.02531645569620253
$ws.Range("P4").Value = 0.759493670886076
$ws.Range("S4").Value = 0.2151898734177215
$ws.Range("J5").Value = 0.3333333333333333
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.07982261640798226
$ws.Range("D6").Value = 0.01330376940133038
$ws.Range("F6").Value = 0.1175166297117517
$ws.Range("J6").Value = 0.1929046563192905
$ws.Range("O6").Value = 0.01995565410199556
$ws.Range("Q6").Value = 0.1906873614190687
$ws.Range("R6").Value = 0.06430155210643015
$ws.Range("S6").Value = 0.3215077605321508
$ws.Range("B7").Value = 0.1114551083591331
$ws.Range("D7").Value = 0.01238390092879257
$ws.Range("F7").Value = 0.06811145510835913
$ws.Range("J7").Value = 0.1238390092879257
$ws.Range("O7").Value = 0.02786377708978328
$ws.Range("Q7").Value = 0.1981424148606811
$ws.Range("R7").Value = 0.108359133126935
$ws.Range("S7").Value = 0.3498452012383901
$ws.Range("B8").Value = 0.1015228426395939
$ws.Range("D8").Value = 0.01649746192893401
$ws.Range("E8").Value = 0.002538071065989848
$ws.Range("F8").Value = 0.07487309644670051
$ws.Range("J8").Value = 0.116751269035533
$ws.Range("O8").Value = 0.03299492385786802
$ws.Range("Q8").Value = 0.1941624365482233
$ws.Range("R8").Value = 0.08756345177664974
$ws.Range("S8").Value = 0.3730964467005076
$ws.Range("B9").Value = 0.1183431952662722
$ws.Range("D9").Value = 0.01183431952662722
$ws.Range("F9").Value = 0.1124260355029586
$ws.Range("J9").Value = 0.06804733727810651
$ws.Range("O9").Value = 0.02662721893491124
$ws.Range("Q9").Value = 0.1538461538461539
$ws.Range("R9").Value = 0.1035502958579882
$ws.Range("S9").Value = 0.4053254437869823
$ws.Range("B10").Value = 0.1212534059945504
$ws.Range("D10").Value = 0.02497729336966394
$ws.Range("E10").Value = 0.0009082652134423251
$ws.Range("F10").Value = 0.07447774750227067
$ws.Range("J10").Value = 0.1049046321525886
$ws.Range("O10").Value = 0.02588555858310627
$ws.Range("Q10").Value = 0.2016348773841962
$ws.Range("R10").Value = 0.07856494096276112
$ws.Range("S10").Value = 0.3673932788374205
$ws.Range("G11").Value = 0.1139489194499018
$ws.Range("J11").Value = 0.106090373280943
$ws.Range("K11").Value = 0.1866404715127701
$ws.Range("L11").Value = 0.5717092337917485
$ws.Range("S11").Value = 0.02161100196463654
$ws.Range("G12").Value = 0.707641196013289
$ws.Range("J12").Value = 0.2159468438538206
$ws.Range("K12").Value = 0.006644518272425249
$ws.Range("L12").Value = 0.03322259136212625
$ws.Range("S12").Value = 0.03654485049833887
$ws.Range("G13").Value = 0.6588235294117647
$ws.Range("J13").Value = 0.3176470588235294
$ws.Range("S13").Value = 0.02352941176470588
$ws.Range("G14").Value = 0.7142857142857143
$ws.Range("J14").Value = 0.1428571428571428
$ws.Range("S14").Value = 0.1428571428571428
$ws.Range("F15").Value = 0.02552204176334107
$ws.Range("H15").Value = 0.1554524361948956
$ws.Range("I15").Value = 0.0765661252900232
$ws.Range("J15").Value = 0.3503480278422274
$ws.Range("K15").Value = 0.06032482598607888
$ws.Range("M15").Value = 0.006960556844547564
$ws.Range("O15").Value = 0.07424593967517401
$ws.Range("S15").Value = 0.2505800464037123
$ws.Range("F16").Value = 0.02439024390243903
$ws.Range("H16").Value = 0.2005420054200542
$ws.Range("I16").Value = 0.08672086720867209
$ws.Range("J16").Value = 0.3983739837398374
$ws.Range("K16").Value = 0.08401084010840108
$ws.Range("M16").Value = 0.02439024390243903
$ws.Range("N16").Value = 0.002710027100271003
$ws.Range("O16").Value = 0.05420054200542006
$ws.Range("S16").Value = 0.1246612466124661
$ws.Range("F17").Value = 0.01884422110552764
$ws.Range("H17").Value = 0.1670854271356784
$ws.Range("I17").Value = 0.09170854271356783
$ws.Range("J17").Value = 0.4422110552763819
$ws.Range("K17").Value = 0.09547738693467336
$ws.Range("M17").Value = 0.01507537688442211
$ws.Range("N17").Value = 0.001256281407035176
$ws.Range("O17").Value = 0.06030150753768844
$ws.Range("S17").Value = 0.1080402010050251
$ws.Range("F18").Value = 0.008797653958944282
$ws.Range("H18").Value = 0.187683284457478
$ws.Range("I18").Value = 0.07038123167155426
$ws.Range("J18").Value = 0.4398826979472141
$ws.Range("K18").Value = 0.1055718475073314
$ws.Range("M18").Value = 0.01466275659824047
$ws.Range("N18").Value = 0.005865102639296188
$ws.Range("O18").Value = 0.05278592375366569
$ws.Range("S18").Value = 0.1143695014662757
$ws.Range("F19").Value = 0.01948051948051948
$ws.Range("H19").Value = 0.2101113172541744
$ws.Range("I19").Value = 0.08256029684601113
$ws.Range("J19").Value = 0.3654916512059369
$ws.Range("K19").Value = 0.1108534322820037
$ws.Range("M19").Value = 0.02782931354359926
$ws.Range("N19").Value = 0.001855287569573284
$ws.Range("O19").Value = 0.06725417439703155
$ws.Range("S19").Value = 0.1145640074211503
